$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 1678
$ws.Range("I2").Value = 4430
$ws.Range("J2").Value = 18686
$ws.Range("K2").Value = 93
$ws.Range("L2").Value = 5171
$ws.Range("M2").Value = 325
$ws.Range("N2").Value = 3217
$ws.Range("O2").Value = 12
$ws.Range("P2").Value = 80
$ws.Range("Q2").Value = 41
$ws.Range("R2").Value = 243
$ws.Range("S2").Value = 1985
$ws.Range("T2").Value = 3312
$ws.Range("U2").Value = 257
$ws.Range("V2").Value = 28971
$ws.Range("W2").Value = 13
$ws.Range("X2").Value = 28970
$ws.Range("Y2").Value = 39
$ws.Range("Z2").Value = 434
$ws.Range("AA2").Value = 191
